# Slide 12, shape "מציין מיקום תוכן 2" (content placeholder) —
# merge the "LSTM + " / "Fasttext" / " embeddings" runs into a single
# run reading "LSTM + Glove embeddings".
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

$oldText = "LSTM + Fasttext embeddings"
$newText = "LSTM + Glove embeddings"

# Select exactly the characters spanning the three existing runs and
# replace them with the new text in one go; this collapses the run
# split (and the err="1" spell-check flag on "Fasttext") into a single
# run that inherits the first run's formatting.
$sub = $tr.Characters(1, $oldText.Length)
$sub.Text = $newText
